$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71: new data for 2020-08-09 (Raw and Clean Data from SSA for August 9th)
# Column A needs to stay a literal text string "2020-08-09" (matching the
# existing Fecha column, which stores dates as shared-string text rather than
# real date serials). Assigning the literal text directly triggers this
# engine's date auto-detection and silently turns it into a date serial, so
# instead we put a text-formula in the cell ( ="2020-08-09" ) and then
# Copy/PasteSpecial-values it into itself: that freezes the formula's cached
# string result as a plain value without ever touching/creating any cell
# style (no stray NumberFormat entries left behind in styles.xml).
$ws.Range("A71").Formula = "=""2020-08-09"""
$ws.Range("A71").Copy()
$ws.Range("A71").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B71").Value = 480278
$ws.Range("C71").Value = 526911
$ws.Range("D71").Value = 84506
$ws.Range("E71").Value = 52298
$ws.Range("F71").Value = 26.61
